$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a paragraph whose visible text begins with a given prefix.
# Using a text anchor (rather than a hard-coded paragraph index) keeps the
# script resilient to the paragraph collection being re-numbered as earlier
# edits insert new paragraphs.
# ---------------------------------------------------------------------------
function Find-ParagraphByPrefix($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $t = $para.Range.Text
        if ($t.Length -ge $prefix.Length -and $t.Substring(0, $prefix.Length) -eq $prefix) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Step 1: "2.2.3. Бібліотеки WiringPi та softPwm" becomes two paragraphs:
#   2.2.3. Система збірки Bazel
#   2.2.4. Бібліотеки WiringPi та softPwm   (renumbered, indented like its
#                                            2.2.2.x siblings)
# We rebuild the paragraph (and the new one that follows it) from raw OOXML
# via Range.InsertXML so every run keeps its original font/size/rsid
# metadata and the language (w:lang) is correct per run.
# ---------------------------------------------------------------------------
$target1 = Find-ParagraphByPrefix $d "2.2.3. Бібліотеки"
if ($null -eq $target1) {
    $target1 = Find-ParagraphByPrefix $d "2.2.3."
}

$xml1 = @'
<w:p w14:paraId="6B2BE537" w14:textId="77BC1CFF" w:rsidR="00AA2E70" w:rsidRPr="00FE41EE" w:rsidRDefault="00AA2E70" w:rsidP="0003499B"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve">2.2.3. Система збірки </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Bazel</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t>2.2.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00270679"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00FE41EE"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve">Бібліотеки </w:t></w:r><w:r w:rsidR="00FE41EE"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>W</w:t></w:r><w:r w:rsidR="00FE41EE" w:rsidRPr="00856549"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>iringPi</w:t></w:r><w:r w:rsidR="00FE41EE" w:rsidRPr="00FE41EE"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00FE41EE"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve">та </w:t></w:r><w:r w:rsidR="00FE41EE" w:rsidRPr="00856549"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t>softPwm</w:t></w:r></w:p>
'@

$target1.Range.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: the page now reflows by one paragraph, so the rendered page break
# that used to fall at the start of "2.4. Висновки до розділу" now falls at
# the start of "2.3. Опис нейронної мережі ..." instead. Move the
# <w:lastRenderedPageBreak/> marker accordingly.
# ---------------------------------------------------------------------------
$target2 = Find-ParagraphByPrefix $d "2.3. Опис нейронної"
if ($null -eq $target2) {
    $target2 = Find-ParagraphByPrefix $d "2.3."
}

$xml2 = @'
<w:p w14:paraId="55E760EE" w14:textId="62A49B47" w:rsidR="001D7AD5" w:rsidRDefault="001D7AD5"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:lastRenderedPageBreak/><w:t>2.3. Опис нейронної мережі – алгоритму зворотнього розповсюдження помилки</w:t></w:r><w:r w:rsidR="00F7687F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@

$target2.Range.InsertXML($xml2) | Out-Null

$target3 = Find-ParagraphByPrefix $d "2.4. Висновки"
if ($null -eq $target3) {
    $target3 = Find-ParagraphByPrefix $d "2.4."
}

$xml3 = @'
<w:p w14:paraId="5FE5E32F" w14:textId="6FA9B851" w:rsidR="00672AD5" w:rsidRDefault="00672AD5"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t>2.4. Висновки до розділу</w:t></w:r><w:r w:rsidR="00B60C3B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="uk-UA"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@

$target3.Range.InsertXML($xml3) | Out-Null

Write-Host "Edit complete."
